# Auto-generated edit script applying the Ultima_Profits.xlsx leve-profit data refresh
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 200
$ws.Range("I8").Value = 200
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 600
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -461
$ws.Range("H100").Value = 2586
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H113").Value = 3549.1667
$ws.Range("I113").Value = 3259
$ws.Range("K113").Value = 3259
$ws.Range("M113").Value = -5
$ws.Range("H116").Value = 2628.5715
$ws.Range("I116").Value = 2125
$ws.Range("J116").Value = 3300
$ws.Range("K116").Value = 2125
$ws.Range("L116").Value = 3300
$ws.Range("M116").Value = 1317
$ws.Range("N116").Value = -10184
$ws.Range("H129").Value = 1539.8474
$ws.Range("I129").Value = 342.92856
$ws.Range("J129").Value = 1912.2222
$ws.Range("K129").Value = 1028.78568
$ws.Range("L129").Value = 5736.6666
$ws.Range("M129").Value = 3971.21432
$ws.Range("N129").Value = -15736.6666
$ws.Range("H137").Value = 2204.9524
$ws.Range("I137").Value = 1045.1111
$ws.Range("J137").Value = 3074.8333
$ws.Range("K137").Value = 3135.3333
$ws.Range("L137").Value = 9224.499899999999
$ws.Range("M137").Value = -585.3333000000002
$ws.Range("N137").Value = -14324.4999
$ws.Range("H138").Value = 1788.8096
$ws.Range("I138").Value = 992.90247
$ws.Range("J138").Value = 3272.0908
$ws.Range("K138").Value = 2978.70741
$ws.Range("L138").Value = 9816.2724
$ws.Range("M138").Value = 2161.29259
$ws.Range("N138").Value = -20096.2724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 20836254
$ws.Range("I132").Value = 41668884
$ws.Range("K132").Value = 125006652
$ws.Range("M132").Value = -125004122

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 517.44446
$ws.Range("I7").Value = 539.25
$ws.Range("K7").Value = 539.25
$ws.Range("M7").Value = -426.25
$ws.Range("H94").Value = 1048
$ws.Range("I94").Value = 920
$ws.Range("J94").Value = 1304
$ws.Range("K94").Value = 920
$ws.Range("L94").Value = 1304
$ws.Range("M94").Value = -469
$ws.Range("N94").Value = -2206
$ws.Range("H99").Value = 1185
$ws.Range("I99").Value = 1185
$ws.Range("K99").Value = 1185
$ws.Range("M99").Value = 313
$ws.Range("H105").Value = 3411.6487
$ws.Range("J105").Value = 4742.857
$ws.Range("L105").Value = 4742.857
$ws.Range("N105").Value = -8236.857
$ws.Range("H107").Value = 2002.75
$ws.Range("I107").Value = 1837
$ws.Range("K107").Value = 1837
$ws.Range("M107").Value = 83

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1909.9166
$ws.Range("I99").Value = 1865.3636
$ws.Range("K99").Value = 1865.3636
$ws.Range("M99").Value = -367.3635999999999
$ws.Range("H107").Value = 742
$ws.Range("I107").Value = 534.7143
$ws.Range("J107").Value = 1225.6666
$ws.Range("K107").Value = 534.7143
$ws.Range("L107").Value = 1225.6666
$ws.Range("M107").Value = 1385.2857
$ws.Range("N107").Value = -5065.6666
$ws.Range("H112").Value = 36702
$ws.Range("J112").Value = 36702
$ws.Range("L112").Value = 36702
$ws.Range("N112").Value = -39656
$ws.Range("H122").Value = 2379.25
$ws.Range("I122").Value = 2379.25
$ws.Range("K122").Value = 7137.75
$ws.Range("M122").Value = -4687.75
$ws.Range("H126").Value = 1909.9166
$ws.Range("I126").Value = 1865.3636
$ws.Range("K126").Value = 5596.0908
$ws.Range("M126").Value = -3126.0908
$ws.Range("H132").Value = 3666.682
$ws.Range("I132").Value = 3568.8235
$ws.Range("J132").Value = 3999.4
$ws.Range("K132").Value = 10706.4705
$ws.Range("L132").Value = 11998.2
$ws.Range("M132").Value = -8176.470499999999
$ws.Range("N132").Value = -17058.2
$ws.Range("H134").Value = 955270.4
$ws.Range("I134").Value = 2939.8096
$ws.Range("J134").Value = 5955006
$ws.Range("K134").Value = 8819.4288
$ws.Range("L134").Value = 17865018
$ws.Range("M134").Value = -6284.4288
$ws.Range("N134").Value = -17870088

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1164.138
$ws.Range("I132").Value = 671.5789
$ws.Range("J132").Value = 2100
$ws.Range("K132").Value = 6044.2101
$ws.Range("L132").Value = 18900
$ws.Range("M132").Value = -3514.2101
$ws.Range("N132").Value = -23960

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2107.25
$ws.Range("I122").Value = 1671
$ws.Range("J122").Value = 2718
$ws.Range("K122").Value = 5013
$ws.Range("L122").Value = 8154
$ws.Range("M122").Value = -2563
$ws.Range("N122").Value = -13054
$ws.Range("H126").Value = 4293.55
$ws.Range("I126").Value = 1287.3334
$ws.Range("K126").Value = 3862.0002
$ws.Range("M126").Value = -1392.0002
$ws.Range("H132").Value = 4051.5454
$ws.Range("I132").Value = 3852.25
$ws.Range("J132").Value = 4583
$ws.Range("K132").Value = 11556.75
$ws.Range("L132").Value = 13749
$ws.Range("M132").Value = -9026.75
$ws.Range("N132").Value = -18809
$ws.Range("H133").Value = 44536.184
$ws.Range("J133").Value = 44536.184
$ws.Range("L133").Value = 44536.184
$ws.Range("N133").Value = -54656.184

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 675.5
$ws.Range("I2").Value = 501
$ws.Range("J2").Value = 850
$ws.Range("K2").Value = 501
$ws.Range("L2").Value = 850
$ws.Range("M2").Value = -389
$ws.Range("N2").Value = -1074
$ws.Range("H7").Value = 5090.6
$ws.Range("I7").Value = 5108.0713
$ws.Range("J7").Value = 5075.3125
$ws.Range("K7").Value = 5108.0713
$ws.Range("L7").Value = 5075.3125
$ws.Range("M7").Value = -4996.0713
$ws.Range("N7").Value = -5299.3125
$ws.Range("H46").Value = 933.5925999999999
$ws.Range("I46").Value = 652
$ws.Range("J46").Value = 997.5909
$ws.Range("K46").Value = 652
$ws.Range("L46").Value = 997.5909
$ws.Range("M46").Value = -464
$ws.Range("N46").Value = -1373.5909
$ws.Range("H122").Value = 6696.9165
$ws.Range("I122").Value = 7468.533
$ws.Range("J122").Value = 5410.8887
$ws.Range("K122").Value = 22405.599
$ws.Range("L122").Value = 16232.6661
$ws.Range("M122").Value = -19955.599
$ws.Range("N122").Value = -21132.6661
$ws.Range("H126").Value = 5090.6
$ws.Range("I126").Value = 5108.0713
$ws.Range("J126").Value = 5075.3125
$ws.Range("K126").Value = 15324.2139
$ws.Range("L126").Value = 15225.9375
$ws.Range("M126").Value = -12854.2139
$ws.Range("N126").Value = -20165.9375
$ws.Range("H132").Value = 41674650
$ws.Range("I132").Value = 8644.556
$ws.Range("J132").Value = 166672670
$ws.Range("K132").Value = 25933.668
$ws.Range("L132").Value = 500018010
$ws.Range("M132").Value = -23403.668
$ws.Range("N132").Value = -500023070

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 962.8570999999999
$ws.Range("I81").Value = 848
$ws.Range("J81").Value = 1250
$ws.Range("K81").Value = 1696
$ws.Range("L81").Value = 2500
$ws.Range("M81").Value = -635
$ws.Range("N81").Value = -4622
$ws.Range("H84").Value = 962.8570999999999
$ws.Range("I84").Value = 848
$ws.Range("J84").Value = 1250
$ws.Range("K84").Value = 8480
$ws.Range("L84").Value = 12500
$ws.Range("M84").Value = -3176
$ws.Range("N84").Value = -23108
$ws.Range("H123").Value = 29000
$ws.Range("J123").Value = 29000
$ws.Range("L123").Value = 29000
$ws.Range("N123").Value = -38800
$ws.Range("H132").Value = 1183.3112
$ws.Range("I132").Value = 1062.2285
$ws.Range("J132").Value = 1607.1
$ws.Range("K132").Value = 3186.6855
$ws.Range("L132").Value = 4821.299999999999
$ws.Range("M132").Value = -656.6854999999996
$ws.Range("N132").Value = -9881.299999999999

